$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rotate the "Financial impact" breakdown rows (8-11): the row that was in
# row 11 moves up to row 8, and rows 8-10 each shift down by one row.
# Column A ("Financial impact" label) is unchanged in all of these rows.

$origB8 = $ws.Range("B8").Value()
$origC8 = $ws.Range("C8").Value()
$origD8 = $ws.Range("D8").Value()

$origB9 = $ws.Range("B9").Value()
$origC9 = $ws.Range("C9").Value()
$origD9 = $ws.Range("D9").Value()

$origB10 = $ws.Range("B10").Value()
$origC10 = $ws.Range("C10").Value()
$origD10 = $ws.Range("D10").Value()

$origB11 = $ws.Range("B11").Value()
$origC11 = $ws.Range("C11").Value()
$origD11 = $ws.Range("D11").Value()

# New row 8 <- old row 11
$ws.Range("B8").Value = $origB11
$ws.Range("C8").Value = $origC11
$ws.Range("D8").Value = $origD11

# New row 9 <- old row 8
$ws.Range("B9").Value = $origB8
$ws.Range("C9").Value = $origC8
$ws.Range("D9").Value = $origD8

# New row 10 <- old row 9
$ws.Range("B10").Value = $origB9
$ws.Range("C10").Value = $origC9
$ws.Range("D10").Value = $origD9

# New row 11 <- old row 10
$ws.Range("B11").Value = $origB10
$ws.Range("C11").Value = $origC10
$ws.Range("D11").Value = $origD10

# Update the saved selection/active cell to match the author's final state.
$ws.Range("I9").Select()
